# Germany Landesliga - base update (02-03-2024 08:34)
# Swap the data (columns B:AC) between a few pairs of rows - the "id"
# column (A) stays put, only the match data that was attached to the
# wrong id moves to its correct row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

Swap-RowData $ws 12 13
Swap-RowData $ws 40 41
Swap-RowData $ws 46 47

# Append the new fixture as row 62, copying the row-61 formatting first
# (bold/border id style in A, date/time number format in E) then filling
# in the values.
$ws.Range("A61:AC61").Copy()
$ws.Range("A62:AC62").PasteSpecial(-4122)

$ws.Range("A62").Value = 60
$ws.Range("B62").Value = 7897140
$ws.Range("C62").Value = "Germany Landesliga"
$ws.Range("D62").Value = "Germany Landesliga"
$ws.Range("E62").Value = 45352.66666666666
$ws.Range("F62").Value = "SC Dsseldorf West"
$ws.Range("G62").Value = "ASV Suchteln"
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 3
$ws.Range("J62").Value = "A"
$ws.Range("K62").Value = 2.15
$ws.Range("L62").Value = 4.2
$ws.Range("M62").Value = 2.4
$ws.Range("N62").Value = 2.15
$ws.Range("O62").Value = 4.2
$ws.Range("P62").Value = 2.4
$ws.Range("Q62").Value = 0
$ws.Range("R62").Value = 1.8
$ws.Range("S62").Value = 2
$ws.Range("T62").Value = 3.75
$ws.Range("U62").Value = 1.85
$ws.Range("V62").Value = 1.95
$ws.Range("W62").Value = -1
$ws.Range("X62").Value = -1
$ws.Range("Y62").Value = 1.4
$ws.Range("Z62").Value = -1
$ws.Range("AA62").Value = 1
$ws.Range("AB62").Value = -1
$ws.Range("AC62").Value = 0.95
